# "prepared for new Specialisations"
# Updates the "umiejetnasci" (specialisations) worksheet:
#  - tweak two descriptive cells (D3, D4)
#  - give D2 the same wrap-text formatting as the rest of the "opis" column
#  - flesh out the still-empty "Luki i kusze" / "Materialy wybuchowe" /
#    "Pistolety maszynowe" rows with their "umiejka1" (C) values, and the
#    "Luki i kusze" row also gets its "opis" (D) value
#  - move the saved cursor/selection back up to the top of the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- text fixes -----------------------------------------------------------
$ws.Range("D3").Value = "Amunicja: pośrednia, karabinowa "
$ws.Range("D4").Value = "Pistolet, Rewolwer, Działo ręczne "

# --- formatting fix: D2 should wrap like the other "opis" cells -----------
$ws.Range("D2").WrapText = $true

# --- row 7 ("Luki i kusze") gains an "umiejka1" + its description ---------
$ws.Rows.Item(7).RowHeight = 23.85
$ws.Range("C7").WrapText = $true
$ws.Range("C7").Value = "Wprawa`n"
$ws.Range("D7").Value = "łuk, Kusza"

# --- row 8 ("Materialy wybuchowe") gains an "umiejka1" ---------------------
$ws.Range("C8").Value = "Wprawa "

# --- row 9 ("Pistolety maszynowe") gains an "umiejka1" ---------------------
$ws.Range("C9").Value = "Wprawa "

# --- restore the view to the top of the sheet with C3 selected ------------
$aw = $excel.ActiveWindow
$aw.ScrollRow = 1
$aw.ScrollColumn = 1
$ws.Range("C3").Select()
